# Saldo_guide.xlsx update: refresh "Dt. Referencia" (column G) from
# 2024-08-09 to 2024-08-12 for every data row, update the handful of
# account balances whose projected/forecast split changed between the
# two extracts, and rename the sheet/tab to match the new extraction
# timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data rows (2 through 274) get the new reference date.
$ws.Range("G2:G274").Value = 45516

# Rows whose "Vl. Projetado" (D) / "Saldo Previsto" (E) / "Vl. Total" (H)
# values were recalculated in the new extract.
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = -2994.88

$ws.Cells.Item(6, 5).Value = 442.95
$ws.Cells.Item(6, 8).Value = 442.95

$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 647.54999999999995

$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = -140.76

$ws.Cells.Item(245, 5).Value = 62.04
$ws.Cells.Item(245, 8).Value = 62.04

$ws.Cells.Item(271, 4).Value = 0
$ws.Cells.Item(271, 5).Value = -1242.02

# Sheet/tab name reflects the newer export timestamp.
$ws.Name = "IClientBalance-20240812-102707-"
